# Apply updated Ligand/Receptor expressing-cell counts and recomputed
# dependent statistics (Natmi analysis update per Dr Hou advice).
$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item(1)

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.747119
$ws.Range("H2").Value = 2.241357
$ws.Range("I2").Value = 0.03096954854571248
$ws.Range("J2").Value = 0.03096954854571248
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.8369776666666665
$ws.Range("N2").Value = 2.510933
$ws.Range("O2").Value = 0.0694586718035551
$ws.Range("P2").Value = 0.06945867180355511
$ws.Range("Q2").Value = 0.6253219173423332
$ws.Range("R2").Value = 5.627897256080999
$ws.Range("S2").Value = 0.00215110370834091
$ws.Range("T2").Value = 0.00215110370834091

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.747119
$ws.Range("H3").Value = 2.241357
$ws.Range("I3").Value = 0.03096954854571248
$ws.Range("J3").Value = 0.03096954854571248
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.427350333333333
$ws.Range("N3").Value = 7.282051
$ws.Range("O3").Value = 0.2014397000898671
$ws.Range("P3").Value = 0.2014397000898671
$ws.Range("Q3").Value = 1.813519553689667
$ws.Range("R3").Value = 16.321675983207
$ws.Range("S3").Value = 0.006238496570966903
$ws.Range("T3").Value = 0.006238496570966903

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.747119
$ws.Range("H4").Value = 2.241357
$ws.Range("I4").Value = 0.03096954854571248
$ws.Range("J4").Value = 0.03096954854571248
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 8.785681666666667
$ws.Range("N4").Value = 26.357045
$ws.Range("O4").Value = 0.7291016281065776
$ws.Range("P4").Value = 0.7291016281065776
$ws.Range("Q4").Value = 6.563949701118333
$ws.Range("R4").Value = 59.07554731006499
$ws.Range("S4").Value = 0.02257994826640466
$ws.Range("T4").Value = 0.02257994826640466

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 19.74619233333334
$ws.Range("H5").Value = 59.23857700000001
$ws.Range("I5").Value = 0.8185184181638298
$ws.Range("J5").Value = 0.8185184181638298
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.8369776666666665
$ws.Range("N5").Value = 2.510933
$ws.Range("O5").Value = 0.0694586718035551
$ws.Range("P5").Value = 0.06945867180355511
$ws.Range("Q5").Value = 16.52712198470455
$ws.Range("R5").Value = 148.744097862341
$ws.Range("S5").Value = 0.05685320217240653
$ws.Range("T5").Value = 0.05685320217240653

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 19.74619233333334
$ws.Range("H6").Value = 59.23857700000001
$ws.Range("I6").Value = 0.8185184181638298
$ws.Range("J6").Value = 0.8185184181638298
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.427350333333333
$ws.Range("N6").Value = 7.282051
$ws.Range("O6").Value = 0.2014397000898671
$ws.Range("P6").Value = 0.2014397000898671
$ws.Range("Q6").Value = 47.93092654238079
$ws.Range("R6").Value = 431.3783388814271
$ws.Range("S6").Value = 0.1648821046729543
$ws.Range("T6").Value = 0.1648821046729543

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 19.74619233333334
$ws.Range("H7").Value = 59.23857700000001
$ws.Range("I7").Value = 0.8185184181638298
$ws.Range("J7").Value = 0.8185184181638298
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 8.785681666666667
$ws.Range("N7").Value = 26.357045
$ws.Range("O7").Value = 0.7291016281065776
$ws.Range("P7").Value = 0.7291016281065776
$ws.Range("Q7").Value = 173.4837599694406
$ws.Range("R7").Value = 1561.353839724965
$ws.Range("S7").Value = 0.5967831113184688
$ws.Range("T7").Value = 0.5967831113184688

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 3.630999
$ws.Range("H8").Value = 10.892997
$ws.Range("I8").Value = 0.1505120332904577
$ws.Range("J8").Value = 0.1505120332904577
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.8369776666666665
$ws.Range("N8").Value = 2.510933
$ws.Range("O8").Value = 0.0694586718035551
$ws.Range("P8").Value = 0.06945867180355511
$ws.Range("Q8").Value = 3.039065070688999
$ws.Range("R8").Value = 27.35158563620099
$ws.Range("S8").Value = 0.01045436592280766
$ws.Range("T8").Value = 0.01045436592280766

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 3.630999
$ws.Range("H9").Value = 10.892997
$ws.Range("I9").Value = 0.1505120332904577
$ws.Range("J9").Value = 0.1505120332904577
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.427350333333333
$ws.Range("N9").Value = 7.282051
$ws.Range("O9").Value = 0.2014397000898671
$ws.Range("P9").Value = 0.2014397000898671
$ws.Range("Q9").Value = 8.813706632983
$ws.Range("R9").Value = 79.32335969684699
$ws.Range("S9").Value = 0.03031909884594589
$ws.Range("T9").Value = 0.0303190988459459

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 3.630999
$ws.Range("H10").Value = 10.892997
$ws.Range("I10").Value = 0.1505120332904577
$ws.Range("J10").Value = 0.1505120332904577
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 8.785681666666667
$ws.Range("N10").Value = 26.357045
$ws.Range("O10").Value = 0.7291016281065776
$ws.Range("P10").Value = 0.7291016281065776
$ws.Range("Q10").Value = 31.900801345985
$ws.Range("R10").Value = 287.107212113865
$ws.Range("S10").Value = 0.1097385685217041
$ws.Range("T10").Value = 0.1097385685217041
